# Weekly data refresh: rows 2-16 (Fecha, Volumen, Precios, Unidad, Origen, etc.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45086
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 18000
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("R2").Value = 'Región del Maule'
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 45099
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 45113
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 13000
$ws.Range("Q4").Value = '$/caja 12 kilos granel'
$ws.Range("R4").Value = 'Región del Maule'
$ws.Range("S4").Value = 13000
$ws.Range("T4").Value = 1

# Row 5
$ws.Range("D5").Value = 45096
$ws.Range("M5").Value = 50
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 18000
$ws.Range("S5").Value = 1000

# Row 6
$ws.Range("D6").Value = 44698
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 16500
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 917

# Row 7
$ws.Range("D7").Value = 45111
$ws.Range("M7").Value = 30
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range("Q7").Value = '$/caja 12 kilos granel'
$ws.Range("R7").Value = 'Región del Maule'
$ws.Range("S7").Value = 13000
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("D8").Value = 45092
$ws.Range("M8").Value = 35
$ws.Range("O8").Value = 19000
$ws.Range("P8").Value = 18571
$ws.Range("S8").Value = 1032

# Row 10
$ws.Range("D10").Value = 44334
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 12500
$ws.Range("Q10").Value = '$/caja 12 kilos empedrada'
$ws.Range("S10").Value = 1042
$ws.Range("T10").Value = 12

# Row 11
$ws.Range("D11").Value = 45085
$ws.Range("M11").Value = 60
$ws.Range("O11").Value = 19000
$ws.Range("P11").Value = 18500
$ws.Range("S11").Value = 1028

# Row 12
$ws.Range("D12").Value = 45107
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("S12").Value = 1000
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44330
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 16000
$ws.Range("P13").Value = 15500
$ws.Range("Q13").Value = '$/caja 18 kilos granel'
$ws.Range("R13").Value = 'Provincia de Curicó'
$ws.Range("S13").Value = 861
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 60
$ws.Range("N14").Value = 17500
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 17750
$ws.Range("S14").Value = 1109

# Row 15
$ws.Range("D15").Value = 44316
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = 16000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 16000
$ws.Range("Q15").Value = '$/caja 16 kilos granel'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("T15").Value = 16

# Row 16
$ws.Range("D16").Value = 44344
$ws.Range("M16").Value = 120
$ws.Range("N16").Value = 13000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 13500
$ws.Range("R16").Value = 'Provincia de Curicó'
$ws.Range("S16").Value = 750

Write-Output "Updated rows 2-16 per weekly refresh."
